$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (shifts un_franzosa_ControlvsDisease_Age and everything below down by 1)
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "un_franzosa_ControlvsCD_ConvCD"
$ws.Range("B9").Value = 0.04
$ws.Range("C9").Value = 0.1
$ws.Range("D9").Value = 0.02
$ws.Range("E9").Value = 0.36
$ws.Range("F9").Value = 0.84
$ws.Range("G9").Value = 0.5
$ws.Range("H9").Value = 0.58

# Insert a new row at row 14 (after un_franzosa_ControlvsUC_Age which is now row 13, before un_franzosa_ControlvsUC_Fp)
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = "un_franzosa_ControlvsUC_ConvUC"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0.04
$ws.Range("E14").Value = 0.3
$ws.Range("F14").Value = 0.96
$ws.Range("G14").Value = 0.7
$ws.Range("H14").Value = 0.66
